# Replace the English "Note from Chuck..." notes body text on slide 1 with
# the Greek translation, as captured by the canonical OOXML diff.
#
# The target paragraph in the diff is rendered from five separate <a:r> runs
# (they differ only in lang="el-GR"/"en-US" and a stray err="1" flag - the
# visible text and formatting, dk2 scheme color, is identical across all of
# them), so concatenating their <a:t> contents reproduces the exact visible
# notes text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$np = $s.NotesPage
$sh = $np.Shapes.Item(1)

# Concatenation of the five <a:r><a:t> runs from the target OOXML, in order.
$run1 = "Σημείωση από τον "
$run2 = " Chuck"
$run3 = ". Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/"
$run4 = "ες"
$run5 = " αναγνώρισης."

$sh.TextFrame.TextRange.Text = $run1 + $run2 + $run3 + $run4 + $run5
